$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Versão": add a new version-history row (row 3) below the existing
# "01.00 - 06/04/2014" row, reusing the same author name.
# ---------------------------------------------------------------------------
$wsVersao = $wb.Worksheets.Item("Versão")
$wsVersao.Range("A3").Value = "02.00 - 09/04/2014"
$wsVersao.Range("B3").Value = "Rodrigo Melo"

# ---------------------------------------------------------------------------
# Sheet "Scripts": add a new test-case row (row 9 / spreadsheet row 10) that
# covers trying to register an already-existing project.
# ---------------------------------------------------------------------------
$wsScripts = $wb.Worksheets.Item("Scripts")

$objetivo = @'
Efetuar o cadastro de projeto sem sucesso. Projeto já existe.
Passo #1 e #2 executados com sucesso;
Deve existir Projeto cadastrado.
'@

$passos = @'
1- Executar o passo #1 acessando a tela de criação de projetos.
2- Preencher os campos obrigatórios com informações de projetos já cadastrados.
3- Acionar a opção de Salvar.
'@

$resultado = @'
O sistema não deverá permitir salvar o novo projeto informando mensagem de que projeto já está cadastrado.
"Já existe Projeto cadastrado com o nome informado."
'@

$wsScripts.Range("A10").Value = 9
$wsScripts.Range("B10").Value = $objetivo
$wsScripts.Range("C10").Value = $passos
$wsScripts.Range("D10").Value = $resultado

# Grow row 10 so the new multi-line content is fully visible (matches the
# taller rows used for the other test-case entries).
$wsScripts.Rows.Item(10).RowHeight = 63.75
